# Apply the "Deploy the implementation guide" metadata refresh:
#  - Update the "Date" value on the Metadata sheet
#  - Update the "Contact" value on the Metadata sheet
#  - Insert a new "Jurisdiction" row (with an empty value) right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) Refresh the Date value (row 8, column B)
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# 2) Refresh the Contact value (row 10, column B)
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# 3) Insert a new row for "Jurisdiction" right after the "Contact" row (row 10),
#    so it becomes the new row 11, pushing everything below down by one.
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "Jurisdiction"
# A lone leading apostrophe forces Excel to store this as an explicit empty
# text value (shared string) instead of simply clearing the cell.
$ws.Range("B11").Value = "'"

# Copy the formatting from the row that got pushed down (now row 12, formerly
# row 11 "Description") onto the freshly inserted row so the style matches
# the rest of the table (style index 2) instead of Excel's blank default,
# and also clears the quote-prefix flag picked up from the apostrophe above.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
